$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3.685507
$ws.Range("H2").Value = 11.056521
$ws.Range("I2").Value = 0.3585631737883472
$ws.Range("J2").Value = 0.3585631737883472
$ws.Range("M2").Value = 18.62071366666667
$ws.Range("N2").Value = 55.862141
$ws.Range("O2").Value = 0.1068221077965325
$ws.Range("P2").Value = 0.1068221077965325
$ws.Range("Q2").Value = 68.62677056349567
$ws.Range("R2").Value = 617.6409350714611
$ws.Range("S2").Value = 0.03830247400228563
$ws.Range("T2").Value = 0.03830247400228563
$ws.Range("G3").Value = 3.685507
$ws.Range("H3").Value = 11.056521
$ws.Range("I3").Value = 0.3585631737883472
$ws.Range("J3").Value = 0.3585631737883472
$ws.Range("O3").Value = 0.7040307798496723
$ws.Range("P3").Value = 0.7040307798496723
$ws.Range("Q3").Value = 452.2973735962056
$ws.Range("R3").Value = 4070.67636236585
$ws.Range("S3").Value = 0.2524395108675837
$ws.Range("T3").Value = 0.2524395108675837
$ws.Range("G4").Value = 3.685507
$ws.Range("H4").Value = 11.056521
$ws.Range("I4").Value = 0.3585631737883472
$ws.Range("J4").Value = 0.3585631737883472
$ws.Range("M4").Value = 32.97121066666667
$ws.Range("N4").Value = 98.91363200000001
$ws.Range("O4").Value = 0.1891471123537951
$ws.Range("P4").Value = 0.1891471123537951
$ws.Range("Q4").Value = 121.5156277104747
$ws.Range("R4").Value = 1093.640649394272
$ws.Range("S4").Value = 0.06782118891847788
$ws.Range("T4").Value = 0.06782118891847787
$ws.Range("I5").Value = 0.009647184430711629
$ws.Range("J5").Value = 0.009647184430711629
$ws.Range("M5").Value = 18.62071366666667
$ws.Range("N5").Value = 55.862141
$ws.Range("O5").Value = 0.1068221077965325
$ws.Range("P5").Value = 0.1068221077965325
$ws.Range("Q5").Value = 1.846411346473
$ws.Range("R5").Value = 16.617702118257
$ws.Range("S5").Value = 0.001030532575190507
$ws.Range("T5").Value = 0.001030532575190507
$ws.Range("I6").Value = 0.009647184430711629
$ws.Range("J6").Value = 0.009647184430711629
$ws.Range("O6").Value = 0.7040307798496723
$ws.Range("P6").Value = 0.7040307798496723
$ws.Range("S6").Value = 0.006791914778107525
$ws.Range("T6").Value = 0.006791914778107525
$ws.Range("I7").Value = 0.009647184430711629
$ws.Range("J7").Value = 0.009647184430711629
$ws.Range("M7").Value = 32.97121066666667
$ws.Range("N7").Value = 98.91363200000001
$ws.Range("O7").Value = 0.1891471123537951
$ws.Range("P7").Value = 0.1891471123537951
$ws.Range("Q7").Value = 3.269392278496
$ws.Range("R7").Value = 29.424530506464
$ws.Range("S7").Value = 0.001824737077413595
$ws.Range("T7").Value = 0.001824737077413595
$ws.Range("G8").Value = 6.493877
$ws.Range("H8").Value = 19.481631
$ws.Range("I8").Value = 0.6317896417809412
$ws.Range("J8").Value = 0.6317896417809411
$ws.Range("M8").Value = 18.62071366666667
$ws.Range("N8").Value = 55.862141
$ws.Range("O8").Value = 0.1068221077965325
$ws.Range("P8").Value = 0.1068221077965325
$ws.Range("Q8").Value = 120.9206242035523
$ws.Range("R8").Value = 1088.285617831971
$ws.Range("S8").Value = 0.06748910121905634
$ws.Range("T8").Value = 0.06748910121905632
$ws.Range("G9").Value = 6.493877
$ws.Range("H9").Value = 19.481631
$ws.Range("I9").Value = 0.6317896417809412
$ws.Range("J9").Value = 0.6317896417809411
$ws.Range("O9").Value = 0.7040307798496723
$ws.Range("P9").Value = 0.7040307798496723
$ws.Range("Q9").Value = 796.9496494123623
$ws.Range("R9").Value = 7172.54684471126
$ws.Range("S9").Value = 0.4447993542039811
$ws.Range("T9").Value = 0.4447993542039811
$ws.Range("G10").Value = 6.493877
$ws.Range("H10").Value = 19.481631
$ws.Range("I10").Value = 0.6317896417809412
$ws.Range("J10").Value = 0.6317896417809411
$ws.Range("M10").Value = 32.97121066666667
$ws.Range("N10").Value = 98.91363200000001
$ws.Range("O10").Value = 0.1891471123537951
$ws.Range("P10").Value = 0.1891471123537951
$ws.Range("Q10").Value = 214.1109866104214
$ws.Range("R10").Value = 1926.998879493792
$ws.Range("S10").Value = 0.1195011863579036
$ws.Range("T10").Value = 0.1195011863579036
